$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the updated LR-pairs table (Sending cluster x Target cluster = FAPs/sCs/ECs)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#          E..T numeric metrics
$data = @(
    @("ECs","Bmp2","Acvr2b","ECs",2,0.6666666666666666,1.445484,4.336452,0.1286708197254238,0.1286708197254238,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,2.225562086516,20.030058778644,0.04596235286387656,0.04596235286387656),
    @("ECs","Bmp2","Acvr2b","FAPs",2,0.6666666666666666,1.445484,4.336452,0.1286708197254238,0.1286708197254238,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,2.100103711876,18.900933406884,0.04337138399364429,0.04337138399364429),
    @("ECs","Bmp2","Acvr2b","sCs",2,0.6666666666666666,1.445484,4.336452,0.1286708197254238,0.1286708197254238,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,1.904757149492,17.142814345428,0.03933708286790299,0.03933708286790299),
    @("FAPs","Bmp2","Acvr2b","ECs",3,1,6.292313,18.876939,0.5601148623429528,0.5601148623429528,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,9.688058290020335,87.192524610183,0.2000779741843962,0.2000779741843962),
    @("FAPs","Bmp2","Acvr2b","FAPs",3,1,6.292313,18.876939,0.5601148623429528,0.5601148623429528,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,9.141927470373666,82.277347233363,0.1887992695396144,0.1887992695396144),
    @("FAPs","Bmp2","Acvr2b","sCs",3,1,6.292313,18.876939,0.5601148623429528,0.5601148623429528,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,8.291567512052334,74.62410760847101,0.1712376186189423,0.1712376186189423),
    @("sCs","Bmp2","Acvr2b","ECs",3,1,3.496172,10.488516,0.3112143179316233,0.3112143179316232,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,5.382935993161334,48.44642393845201,0.1111685021327148,0.1111685021327147),
    @("sCs","Bmp2","Acvr2b","FAPs",3,1,3.496172,10.488516,0.3112143179316233,0.3112143179316232,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,5.079491571374667,45.715424142372,0.1049017618457398,0.1049017618457398),
    @("sCs","Bmp2","Acvr2b","sCs",3,1,3.496172,10.488516,0.3112143179316233,0.3112143179316232,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,4.607009564169333,41.463086077524,0.09514405395316873,0.09514405395316872)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
